$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 118826.46
$ws.Range("J40").Value = 3904.5454
$ws.Range("L40").Value = 3904.5454
$ws.Range("N40").Value = -4254.5454
$ws.Range("H74").Value = 7448.615
$ws.Range("I74").Value = 5479
$ws.Range("K74").Value = 5479
$ws.Range("M74").Value = -4543
$ws.Range("H77").Value = 7448.615
$ws.Range("I77").Value = 5479
$ws.Range("K77").Value = 27395
$ws.Range("M77").Value = -22715
$ws.Range("H98").Value = 1587.0646
$ws.Range("I98").Value = 1633.3
$ws.Range("K98").Value = 1633.3
$ws.Range("M98").Value = -135.3
$ws.Range("H122").Value = 1587.0646
$ws.Range("I122").Value = 1633.3
$ws.Range("K122").Value = 4899.9
$ws.Range("M122").Value = -2449.9
$ws.Range("H125").Value = 3803.25
$ws.Range("I125").Value = 1666
$ws.Range("J125").Value = 5085.6
$ws.Range("K125").Value = 14994
$ws.Range("L125").Value = 45770.4
$ws.Range("M125").Value = -12534
$ws.Range("N125").Value = -50690.4
$ws.Range("H132").Value = 2085.0356
$ws.Range("I132").Value = 1877.68
$ws.Range("K132").Value = 5633.04
$ws.Range("M132").Value = -3103.04
$ws.Range("H138").Value = 3463.6262
$ws.Range("I138").Value = 2391.8572
$ws.Range("J138").Value = 3640.1528
$ws.Range("K138").Value = 7175.571599999999
$ws.Range("L138").Value = 10920.4584
$ws.Range("M138").Value = -2035.571599999999
$ws.Range("N138").Value = -21200.4584

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 1500
$ws.Range("I11").Value = 1500
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 1500
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -1356
$ws.Range("N11").Value = $null
$ws.Range("H122").Value = 2917.7812
$ws.Range("J122").Value = 4341
$ws.Range("L122").Value = 13023
$ws.Range("N122").Value = -17923
$ws.Range("H132").Value = 2574.3635
$ws.Range("I132").Value = 2417.8708
$ws.Range("K132").Value = 7253.6124
$ws.Range("M132").Value = -4723.6124

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1801.08
$ws.Range("I20").Value = 1899.5883
$ws.Range("J20").Value = 1591.75
$ws.Range("K20").Value = 1899.5883
$ws.Range("L20").Value = 1591.75
$ws.Range("M20").Value = -1652.5883
$ws.Range("N20").Value = -2085.75
$ws.Range("H86").Value = 1832.2142
$ws.Range("I86").Value = 1757.625
$ws.Range("J86").Value = 1931.6666
$ws.Range("K86").Value = 1757.625
$ws.Range("L86").Value = 1931.6666
$ws.Range("M86").Value = -634.625
$ws.Range("N86").Value = -4177.6666
$ws.Range("H89").Value = 1832.2142
$ws.Range("I89").Value = 1757.625
$ws.Range("J89").Value = 1931.6666
$ws.Range("K89").Value = 8788.125
$ws.Range("L89").Value = 9658.333000000001
$ws.Range("M89").Value = -3172.125
$ws.Range("N89").Value = -20890.333
$ws.Range("H94").Value = 1500
$ws.Range("I94").Value = 1500
$ws.Range("K94").Value = 1500
$ws.Range("M94").Value = -1049
$ws.Range("H99").Value = 3160.9285
$ws.Range("I99").Value = 3205.4
$ws.Range("K99").Value = 3205.4
$ws.Range("M99").Value = -1707.4
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").Value = $null
$ws.Range("H140").Value = 118900
$ws.Range("J140").Value = 118900
$ws.Range("L140").Value = 118900
$ws.Range("N140").Value = -129260

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 4373.25
$ws.Range("I105").Value = 5996.5
$ws.Range("K105").Value = 5996.5
$ws.Range("M105").Value = -4249.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 1492.091
$ws.Range("I18").Value = 702
$ws.Range("J18").Value = 2874.75
$ws.Range("K18").Value = 2106
$ws.Range("L18").Value = 8624.25
$ws.Range("M18").Value = -1937
$ws.Range("N18").Value = -8962.25
$ws.Range("H38").Value = 208.4
$ws.Range("I38").Value = 14.5
$ws.Range("J38").Value = 337.66666
$ws.Range("K38").Value = 43.5
$ws.Range("L38").Value = 1012.99998
$ws.Range("M38").Value = 303.5
$ws.Range("N38").Value = -1706.99998
$ws.Range("H68").Value = 2231.2632
$ws.Range("J68").Value = 2708.7273
$ws.Range("L68").Value = 8126.1819
$ws.Range("N68").Value = -9748.1819
$ws.Range("H71").Value = 2231.2632
$ws.Range("J71").Value = 2708.7273
$ws.Range("L71").Value = 24378.5457
$ws.Range("N71").Value = -32490.5457
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").Value = $null
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").Value = $null
$ws.Range("H115").Value = 3209.3333
$ws.Range("I115").Value = 3209.3333
$ws.Range("K115").Value = 9627.999899999999
$ws.Range("M115").Value = -8452.999899999999
$ws.Range("H122").Value = 985.1429000000001
$ws.Range("I122").Value = 445
$ws.Range("J122").Value = 1075.1666
$ws.Range("K122").Value = 4005
$ws.Range("L122").Value = 9676.499400000001
$ws.Range("M122").Value = -1555
$ws.Range("N122").Value = -14576.4994
$ws.Range("H129").Value = 1994.6765
$ws.Range("J129").Value = 2064.3438
$ws.Range("L129").Value = 6193.0314
$ws.Range("N129").Value = -16193.0314
$ws.Range("H134").Value = 2529.8572
$ws.Range("I134").Value = 2529.8572
$ws.Range("K134").Value = 7589.571599999999
$ws.Range("M134").Value = -2519.571599999999
$ws.Range("H140").Value = 1337.6666
$ws.Range("I140").Value = 1337.6666
$ws.Range("K140").Value = 4012.9998
$ws.Range("M140").Value = 1167.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1038.8889
$ws.Range("I97").Value = 1112.3125
$ws.Range("K97").Value = 1112.3125
$ws.Range("M97").Value = -616.3125
$ws.Range("H122").Value = 1423.909
$ws.Range("I122").Value = 1235.6364
$ws.Range("K122").Value = 3706.9092
$ws.Range("M122").Value = -1256.9092

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3000.8462
$ws.Range("I7").Value = 2928.3635
$ws.Range("J7").Value = 3399.5
$ws.Range("K7").Value = 2928.3635
$ws.Range("L7").Value = 3399.5
$ws.Range("M7").Value = -2816.3635
$ws.Range("N7").Value = -3623.5
$ws.Range("H22").Value = 2158.9048
$ws.Range("I22").Value = 1960.3
$ws.Range("K22").Value = 1960.3
$ws.Range("M22").Value = -1665.3
$ws.Range("H27").Value = 2158.9048
$ws.Range("I27").Value = 1960.3
$ws.Range("K27").Value = 1960.3
$ws.Range("M27").Value = -1853.3
$ws.Range("H36").Value = 132000
$ws.Range("J36").Value = 132000
$ws.Range("L36").Value = 132000
$ws.Range("N36").Value = -133124
$ws.Range("H40").Value = 2500
$ws.Range("I40").Value = 2500
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 2500
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -2364
$ws.Range("N40").Value = $null
$ws.Range("H99").Value = 74755.14
$ws.Range("I99").Value = 61867.4
$ws.Range("K99").Value = 61867.4
$ws.Range("M99").Value = -58872.4
$ws.Range("H112").Value = 91541.336
$ws.Range("J112").Value = 91541.336
$ws.Range("L112").Value = 91541.336
$ws.Range("N112").Value = -94495.336
$ws.Range("H122").Value = 8210.333000000001
$ws.Range("I122").Value = 5741
$ws.Range("K122").Value = 17223
$ws.Range("M122").Value = -14773
$ws.Range("H124").Value = 108995
$ws.Range("I124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("M124").Value = $null
$ws.Range("H126").Value = 3000.8462
$ws.Range("I126").Value = 2928.3635
$ws.Range("J126").Value = 3399.5
$ws.Range("K126").Value = 8785.0905
$ws.Range("L126").Value = 10198.5
$ws.Range("M126").Value = -6315.0905
$ws.Range("N126").Value = -15138.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 10098
$ws.Range("I96").Value = 4653.75
$ws.Range("K96").Value = 4653.75
$ws.Range("M96").Value = -3280.75
$ws.Range("H126").Value = 4171.684
$ws.Range("I126").Value = 4125.6665
$ws.Range("K126").Value = 12376.9995
$ws.Range("M126").Value = -9906.999500000002
$ws.Range("H132").Value = 2948.8928
$ws.Range("I132").Value = 2798.6667
$ws.Range("K132").Value = 8396.000100000001
$ws.Range("M132").Value = -5866.000100000001
